$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").Value = 7
$ws.Range("E5").Value = 7
$ws.Range("F5").Value = 7

$ws.Range("D6").Value = 8
$ws.Range("E6").Value = 8
$ws.Range("F6").Value = 8
